$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp
$ws.Range("A1").Value = "Datos actualizados a 28 de Junio de 2020 a las 10:12"

# Refresh per-country case counts from the latest snapshot.
# Several countries changed rank (and therefore row position, since rows
# stay sorted descending by total cases) after the update.
function Set-CountryRow($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
  $data = New-Object 'object[,]' 1,8
  $data[0,0] = $country
  $data[0,1] = $total
  $data[0,2] = $nuevos
  $data[0,3] = $activos
  $data[0,4] = $recuperados
  $data[0,5] = $criticos
  $data[0,6] = $muertesHoy
  $data[0,7] = $muertes
  $ws.Range("A" + $row + ":H" + $row).Value = $data
}

Set-CountryRow 6 "Rusia" 634437 6791 399087 226277 0 104 9073
Set-CountryRow 7 "India" 529889 312 310236 203541 0 9 16112
Set-CountryRow 36 "Singapur" 43459 213 37163 6270 0 0 26
Set-CountryRow 37 "Irak" 43262 0 19938 21664 0 0 1660
Set-CountryRow 38 "Ucrania" 42982 917 18934 22919 0 19 1129
Set-CountryRow 42 "Polonia" 33714 0 20548 11731 0 0 1435
Set-CountryRow 44 "Afganistan" 30967 351 12588 17642 0 34 737
Set-CountryRow 53 "Israel" 23497 76 17019 6160 0 1 318
Set-CountryRow 54 "Kazajistan" 20780 461 12824 7783 0 7 173
Set-CountryRow 79 "El Salvador" 5934 207 3557 2225 0 9 152
Set-CountryRow 80 "Republica de Macedonia" 5906 0 2236 3393 0 0 277
Set-CountryRow 81 "Kenia" 5811 0 1936 3734 0 0 141
Set-CountryRow 82 "Tayikistan" 5799 0 4391 1356 0 0 52
Set-CountryRow 83 "Haiti" 5777 55 706 4971 0 2 100
Set-CountryRow 92 "Hungria" 4142 4 2685 876 0 3 581
Set-CountryRow 111 "Estonia" 1987 1 1818 100 0 0 69
Set-CountryRow 114 "Estado de Palestina" 1854 39 447 1403 0 0 4
Set-CountryRow 115 "Islandia" 1836 0 1814 12 0 0 10
Set-CountryRow 116 "Lituania" 1815 2 1503 234 0 0 78
Set-CountryRow 118 "Eslovaquia" 1664 7 1461 175 0 0 28
Set-CountryRow 140 "Uganda" 859 11 794 65 0 0 0
Set-CountryRow 141 "Principado de Andorra" 855 0 799 4 0 0 52
Set-CountryRow 202 "Laos" 19 0 19 0 0 0 0
Set-CountryRow 203 "Santa Lucia" 19 0 19 0 0 0 0
Set-CountryRow 209 "Groenlandia" 13 0 13 0 0 0 0
Set-CountryRow 210 "Islas Malvinas" 13 0 13 0 0 0 0

Write-Host "done"